# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Home win": 2 -> 3 data rows. Reorder existing rows and append a new
# match row (Israel - State Cup).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home win")

$ws.Range("A2").Value = "29-12-2024 12:30"
$ws.Range("B2").Value = "ENGLAND"
$ws.Range("C2").Value = "LEAGUE ONE"
$ws.Range("D2").Value = "Rotherham - Stockport County"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 2.8

$ws.Range("A3").Value = "29-12-2024 15:00"
$ws.Range("B3").Value = "ENGLAND"
$ws.Range("C3").Value = "LEAGUE TWO"
$ws.Range("D3").Value = "Bradford - Chesterfield"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 2.2

$ws.Range("A4").Value = "29-12-2024 18:15"
$ws.Range("B4").Value = "ISRAEL"
$ws.Range("C4").Value = "STATE CUP"
$ws.Range("D4").Value = "Kiryat Yam SC - Beitar Jerusalem"
$ws.Range("E4").Value = 70
$ws.Range("F4").Value = 8.5

# ---------------------------------------------------------------------------
# Sheet "Draw": 4 -> 2 data rows. Drop the two oldest matches, keep and
# update the remaining two.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Draw")

$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

$ws.Range("A2").Value = "29-12-2024 14:00"
$ws.Range("B2").Value = "ITALY"
$ws.Range("C2").Value = "SERIE B"
$ws.Range("D2").Value = "Bari - Spezia"
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 2.85

$ws.Range("A3").Value = "29-12-2024 12:30"
$ws.Range("B3").Value = "CAMEROON"
$ws.Range("C3").Value = "ELITE ONE"
$ws.Range("D3").Value = "Bamboutos - Fauve Azur Elite"
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 2.75

# ---------------------------------------------------------------------------
# Sheet "Btts": 7 -> 3 data rows.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Btts")

$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

$ws.Range("A2").Value = "29-12-2024 15:00"
$ws.Range("B2").Value = "ENGLAND"
$ws.Range("C2").Value = "LEAGUE ONE"
$ws.Range("D2").Value = "Exeter City - Crawley Town"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 1.8

$ws.Range("A3").Value = "29-12-2024 15:00"
$ws.Range("B3").Value = "SCOTLAND"
$ws.Range("C3").Value = "PREMIERSHIP"
$ws.Range("D3").Value = "Motherwell - Rangers"
$ws.Range("E3").Value = 76.7
$ws.Range("F3").Value = 1.95

$ws.Range("A4").Value = "29-12-2024 15:30"
$ws.Range("B4").Value = "PORTUGAL"
$ws.Range("C4").Value = "PRIMEIRA LIGA"
$ws.Range("D4").Value = "Rio Ave - Nacional"
$ws.Range("E4").Value = 86.7
$ws.Range("F4").Value = 1.8

# ---------------------------------------------------------------------------
# Sheet "Over_Under": 5 -> 3 data rows.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Over_Under")

$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

$ws.Range("A2").Value = "29-12-2024 15:00"
$ws.Range("B2").Value = "ENGLAND"
$ws.Range("C2").Value = "CHAMPIONSHIP"
$ws.Range("D2").Value = "Oxford United - Plymouth"
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 1.8
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 3.2

$ws.Range("A3").Value = "29-12-2024 15:00"
$ws.Range("B3").Value = "ENGLAND"
$ws.Range("C3").Value = "LEAGUE ONE"
$ws.Range("D3").Value = "Exeter City - Crawley Town"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 1.85
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 3.25

$ws.Range("A4").Value = "30-12-2024 19:45"
$ws.Range("B4").Value = "NORTHERN-IRELAND"
$ws.Range("C4").Value = "PREMIERSHIP"
$ws.Range("D4").Value = "Loughgall - Coleraine FC"
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 1.73
$ws.Range("G4").Value = 62.5
$ws.Range("H4").Value = 2.88

# ---------------------------------------------------------------------------
# Sheet "Away Win": unchanged by this update.
# ---------------------------------------------------------------------------
